$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new client account (MANOLO 60216) to the PREMIUM group's client list in B3
$cell = $ws.Range("B3")
$cell.Value = $cell.Value2 + ".60216"

# Move the active selection to B4 (as left by the editor after the edit)
$ws.Range("B4").Select()
